# Apply the scheduled-runner profit-recalc updates to the Leve profit sheets.
# Each worksheet (one per crafting class) gets refreshed currentAveragePrice /
# LevePrice / LeveProfit figures (columns H-N) for the rows the runner recomputed.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 90910770
$ws.Range("I100").Value = 250000530
$ws.Range("J100").Value = 2330.2856
$ws.Range("K100").Value = 250000530
$ws.Range("L100").Value = 2330.2856
$ws.Range("M100").Value = -249999989
$ws.Range("N100").Value = -3412.2856
$ws.Range("H116").Value = 4673.1816
$ws.Range("I116").Value = 1200
$ws.Range("K116").Value = 1200
$ws.Range("M116").Value = 2242
$ws.Range("H129").Value = 167762.64
$ws.Range("J129").Value = 189854.27
$ws.Range("L129").Value = 569562.8099999999
$ws.Range("N129").Value = -579562.8099999999
$ws.Range("H138").Value = 1518.0422
$ws.Range("I138").Value = 610.32434
$ws.Range("J138").Value = 2505.853
$ws.Range("K138").Value = 1830.97302
$ws.Range("L138").Value = 7517.559
$ws.Range("M138").Value = 3309.02698
$ws.Range("N138").Value = -17797.559

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4552.72
$ws.Range("I32").Value = 4955.186
$ws.Range("J32").Value = 2080.4285
$ws.Range("K32").Value = 4955.186
$ws.Range("L32").Value = 2080.4285
$ws.Range("M32").Value = -4668.186
$ws.Range("N32").Value = -2654.4285
$ws.Range("H74").Value = 30305250
$ws.Range("I74").Value = 38463656
$ws.Range("J74").Value = 2605.7144
$ws.Range("K74").Value = 38463656
$ws.Range("L74").Value = 2605.7144
$ws.Range("M74").Value = -38462782
$ws.Range("N74").Value = -4353.7144
$ws.Range("H77").Value = 30305250
$ws.Range("I77").Value = 38463656
$ws.Range("J77").Value = 2605.7144
$ws.Range("K77").Value = 192318280
$ws.Range("L77").Value = 13028.572
$ws.Range("M77").Value = -192313912
$ws.Range("N77").Value = -21764.572
$ws.Range("H122").Value = 4116.5625
$ws.Range("I122").Value = 3676.2144
$ws.Range("K122").Value = 11028.6432
$ws.Range("M122").Value = -8578.643199999999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3108
$ws.Range("I20").Value = 3740.182
$ws.Range("K20").Value = 3740.182
$ws.Range("M20").Value = -3493.182
$ws.Range("H107").Value = 1686.6666
$ws.Range("J107").Value = 2110.25
$ws.Range("L107").Value = 2110.25
$ws.Range("N107").Value = -5950.25

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 40681.418
$ws.Range("J20").Value = 40681.418
$ws.Range("L20").Value = 40681.418
$ws.Range("N20").Value = -41153.418
$ws.Range("H30").Value = 40681.418
$ws.Range("J30").Value = 40681.418
$ws.Range("L30").Value = 40681.418
$ws.Range("N30").Value = -40863.418
$ws.Range("H94").Value = 8382.8
$ws.Range("I94").Value = 3200
$ws.Range("K94").Value = 3200
$ws.Range("M94").Value = -2749
$ws.Range("H122").Value = 1399.7307
$ws.Range("I122").Value = 1631.9333
$ws.Range("J122").Value = 1083.091
$ws.Range("K122").Value = 4895.7999
$ws.Range("L122").Value = 3249.273
$ws.Range("M122").Value = -2445.7999
$ws.Range("N122").Value = -8149.272999999999
$ws.Range("H128").Value = 40681.418
$ws.Range("J128").Value = 40681.418
$ws.Range("L128").Value = 40681.418
$ws.Range("N128").Value = -50641.418
$ws.Range("H132").Value = 14382.439
$ws.Range("I132").Value = 17504.129
$ws.Range("J132").Value = 4705.2
$ws.Range("K132").Value = 52512.387
$ws.Range("L132").Value = 14115.6
$ws.Range("M132").Value = -49982.387
$ws.Range("N132").Value = -19175.6
$ws.Range("H134").Value = 909.96
$ws.Range("I134").Value = 797.5714
$ws.Range("K134").Value = 2392.7142
$ws.Range("M134").Value = 142.2857999999997

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 492.18182
$ws.Range("I122").Value = 236.73334
$ws.Range("J122").Value = 1039.5714
$ws.Range("K122").Value = 2130.60006
$ws.Range("L122").Value = 9356.142600000001
$ws.Range("M122").Value = 319.3999400000002
$ws.Range("N122").Value = -14256.1426
$ws.Range("H129").Value = 264063.06
$ws.Range("J129").Value = 313429.25
$ws.Range("L129").Value = 940287.75
$ws.Range("N129").Value = -950287.75
$ws.Range("H131").Value = 808.38
$ws.Range("I131").Value = 614.75
$ws.Range("J131").Value = 816.44794
$ws.Range("K131").Value = 1844.25
$ws.Range("L131").Value = 2449.34382
$ws.Range("M131").Value = 3195.75
$ws.Range("N131").Value = -12529.34382

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 25001132
$ws.Range("I102").Value = 29412828
$ws.Range("K102").Value = 29412828
$ws.Range("M102").Value = -29411206
$ws.Range("H113").Value = 2377.3684
$ws.Range("I113").Value = 2227.05
$ws.Range("J113").Value = 2544.389
$ws.Range("K113").Value = 2227.05
$ws.Range("L113").Value = 2544.389
$ws.Range("M113").Value = -57.05000000000018
$ws.Range("N113").Value = -6884.389
$ws.Range("H122").Value = 60608044
$ws.Range("I122").Value = 41667850
$ws.Range("K122").Value = 125003550
$ws.Range("M122").Value = -125001100
$ws.Range("H126").Value = 5392.56
$ws.Range("I126").Value = 4326.6665
$ws.Range("K126").Value = 12979.9995
$ws.Range("M126").Value = -10509.9995
$ws.Range("H132").Value = 26647.092
$ws.Range("I132").Value = 4595
$ws.Range("K132").Value = 13785
$ws.Range("M132").Value = -11255

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 503.07693
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H55").Value = 60.206898
$ws.Range("I55").Value = 53.46154
$ws.Range("J55").Value = 65.6875
$ws.Range("K55").Value = 53.46154
$ws.Range("L55").Value = 65.6875
$ws.Range("M55").Value = 119.53846
$ws.Range("N55").Value = -411.6875
$ws.Range("H122").Value = 678696.3
$ws.Range("I122").Value = 1402979.2
$ws.Range("J122").Value = 2698.8667
$ws.Range("K122").Value = 4208937.6
$ws.Range("L122").Value = 8096.6001
$ws.Range("M122").Value = -4206487.6
$ws.Range("N122").Value = -12996.6001
$ws.Range("H132").Value = 2270.8235
$ws.Range("I132").Value = 1510.6
$ws.Range("J132").Value = 3356.8572
$ws.Range("K132").Value = 4531.799999999999
$ws.Range("L132").Value = 10070.5716
$ws.Range("M132").Value = -2001.799999999999
$ws.Range("N132").Value = -15130.5716

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H132").Value = 1721.1052
$ws.Range("I132").Value = 1060.3
$ws.Range("J132").Value = 2110.25
$ws.Range("K132").Value = 3180.9
$ws.Range("M132").Value = -650.8999999999996
